$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "difficulty" numbers in column F for rows 9-15
$ws.Range("F9").Value = 2
$ws.Range("F10").Value = 1.2
$ws.Range("F11").Value = 2.3
$ws.Range("F12").Value = 1.2
$ws.Range("F13").Value = 2
$ws.Range("F14").Value = 1.2
$ws.Range("F15").Value = 1.2

# Update the comment text in column G to match the new remarks
$ws.Range("G9").Value = "podobne 1, ale kratsi"
$ws.Range("G11").Value = "da se v tom chybovat ale jde docela"
$ws.Range("G13").Value = "podobné 11"
$ws.Range("G14").Value = "hodne v pohodě"

# Restore the view/selection state captured when the workbook was saved
$ws.Range("G16").Select()
